# GanttChart8_26.xlsx edit
# - Update plan/actual duration + percent-complete figures for several
#   activities on the "Project Planner" sheet.
# - Add a single-space label in A24 (row for "Activity 20").
# - Refresh the view (zoom/scroll/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")
$ws.Activate() | Out-Null

# Row 5: Research
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 0.6

# Row 6: Preliminary Documentation
$ws.Range("G6").Value = 1

# Row 7: Coding Tutorial
$ws.Range("D7").Value = 7
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 0.7

# Row 8: Top Level
$ws.Range("G8").Value = 1

# Row 9: CPU Core + 3 Stages
$ws.Range("D9").Value = 5
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 0.9

# Row 10: Input/Output
$ws.Range("G10").Value = 0.4

# Row 11: Software and Toolchain Setup
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 5
$ws.Range("G11").Value = 0.8

# Row 12: Modifying Tutorial
$ws.Range("C12").Value = 8
$ws.Range("E12").Value = 8

# Row 13: Writing C Programs
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 6

# Row 16: CPU Documentation
$ws.Range("G16").Value = 0.1

# Row 18: Proposal Writing
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 1

# Row 19: Opcode Fetch
$ws.Range("G19").Value = 1

# Row 20: Opcode Decode
$ws.Range("G20").Value = 1

# Row 21: Opcode Execute
$ws.Range("G21").Value = 1

# Row 22: Testing and Simulation
$ws.Range("C22").Value = 6
$ws.Range("E22").Value = 6

# Row 24 (Activity 20): stray single-space note in column A
$ws.Range("A24").Value = " "

# Refresh window view to match the saved state: zoom + scroll + selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Zoom = 63
$ws.Range("BY8").Select() | Out-Null
